# The TFM_INS scenario sheet has a single data row (row 6) describing a
# FLO_SHAR constraint. The LimType in column G was changed from "FX"
# (fixed) to "LO" (lower bound).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INS_1")
$ws.Range("G6").Value = "LO"
